$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B21").Value = "The UPA Manual of Championship Series Tournament Formats"
$ws.Range("C21").Value = "Eric Simon"
$ws.Range("D21").Value = 2008
$ws.Range("E21").Value = "https://usaultimate.org/wp-content/uploads/2020/11/USAU_TournamentFormats.pdf"
$ws.Range("F21").Value = "Semibrackets"
$ws.Range("G21").Value = "x"
$ws.Range("H21").Value = "ultimate formats"

$ws.Range("F22").Select()
